$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Charge Count Quantity" row above the "Charge Description" row ---
# Row 37 already carries the exact style pattern we need (A16/B20/C16/D16/E14, ht=42),
# so copy it as the donor for the new row 38, then overwrite its contents.
$ws.Rows.Item(37).Copy()
$ws.Rows.Item(38).Insert()

# --- Insert "Statute Section ID" / "Statute or Ordinance Section Number" row
#     right after the "Statute ID" / "Statute or Ordinance Number" row ---
# Row 40 (old row 39, now shifted down by the insert above) carries the exact
# style pattern we need (A16/B15/C16/D16/E14, ht=56) for the new row 41.
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(41).Insert()

$ws.Rows.Item(38).RowHeight = 42
$ws.Rows.Item(41).RowHeight = 56

$ws.Range("B38").Value = "Charge Count Quantity"
$ws.Range("C38").Value = ""
$ws.Range("E38").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Charge[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityChargeAssociation/j:Charge/@structures:ref]/j:ChargeCountQuantity"

# The mapping text for the "Statute ID" row moves to the new Statute Section ID row;
# the Statute ID row itself gets a new, more specific mapping path.
$ws.Range("E40").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Charge[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityChargeAssociation/j:Charge/@structures:ref]/j:ChargeStatute/j:StatuteCodeIdentification/nc:IdentificationID"

$ws.Range("B41").Value = "Statute or Ordinance Section Number"
$ws.Range("C41").Value = "Statute Section ID"
$ws.Range("E41").Value = "/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:Charge[@structures:id=/cq-res-doc:CustodyQueryResults/cq-res-ext:Custody/j:ActivityChargeAssociation/j:Charge/@structures:ref]/j:ChargeStatute/j:StatuteCodeSectionIdentification/nc:IdentificationID"

# --- Remove the two now-redundant blank filler rows that used to sit at 58/59 ---
# (they now sit at 60/60 after the two inserts above shifted everything down by 2)
$ws.Rows.Item(60).Delete()
$ws.Rows.Item(60).Delete()

# --- Toggle the hidden blank-spacer block down by two rows ---
$ws.Rows.Item(61).Hidden = $true
$ws.Rows.Item(273).Hidden = $false
$ws.Rows.Item(274).Hidden = $false
$ws.Rows.Item(287).Hidden = $true
$ws.Rows.Item(288).Hidden = $true

# --- Append two more blank filler rows at the very end of the sheet (384 -> 386) ---
$ws.Rows.Item(383).Copy()
$ws.Rows.Item(384).Insert()
$ws.Rows.Item(383).Copy()
$ws.Rows.Item(384).Insert()

# --- Update the view: scroll/selection now centers on the newly inserted row ---
$ws.Range("C38").Select()
